# Update "想去人数" (F column) values on the 展览 and 全部类型 sheets
# to match the newly generated data snapshot.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6816
$ws1.Range("F3").Value = 0
$ws1.Range("F6").Value = 147
$ws1.Range("F7").Value = 0
$ws1.Range("F8").Value = 0
$ws1.Range("F14").Value = 130
$ws1.Range("F15").Value = 0
$ws1.Range("F16").Value = 375
$ws1.Range("F17").Value = 0
$ws1.Range("F19").Value = 4803
$ws1.Range("F20").Value = 86
$ws1.Range("F21").Value = 70
$ws1.Range("F22").Value = 270
$ws1.Range("F23").Value = 0
$ws1.Range("F24").Value = 0

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 0
$ws4.Range("F6").Value = 147
$ws4.Range("F7").Value = 6454
$ws4.Range("F8").Value = 0
$ws4.Range("F9").Value = 0
$ws4.Range("F12").Value = 0
$ws4.Range("F14").Value = 130
$ws4.Range("F16").Value = 0
$ws4.Range("F19").Value = 0
$ws4.Range("F21").Value = 0
$ws4.Range("F22").Value = 70
$ws4.Range("F23").Value = 270
$ws4.Range("F25").Value = 136

$wb.Save()
